$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Update the byline date: "24 Jun 2020" -> "25 Jun 2020"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("24 Jun 2020", $true, $false, $false, $false, $false,
                         $true, 1, $false, "25 Jun 2020", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Update the git revision note: "b1ce88f" -> "5e87c12"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("b1ce88f", $true, $true, $false, $false, $false,
                         $true, 1, $false, "5e87c12", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Rewrite the covariate-adjustment discussion. The two paragraphs
#    beginning "While we randomized clusters..." and "We followed the
#    intention-to-treat principle..." are collapsed into a single
#    paragraph with new wording for the first part; the sentence about
#    the intention-to-treat principle onward is kept verbatim. Including
#    the paragraph mark ("^p") in the search text merges the two
#    paragraphs into one when the replacement text omits it.
# ---------------------------------------------------------------------------
$oldText = "While we randomized clusters by stratifying on district for " + `
    "administrative reasons, we judged prior to analysis that district may " + `
    "be prognostic for all outcomes, and therefore adjusted for district as " + `
    "a fixed effect in all analyses (CHMP 2015). Because we constrained the " + `
    "randomization by lab availability, clinic size (number of new " + `
    "enrollments), proportion of women aged >40 years, and proportion of " + `
    "primiparous women, we also adjusted for these prognostic variables as " + `
    "fixed effects using individual- rather than cluster-level data where " + `
    "possible (Li 2017).^pWe followed the intention-to-treat principle for " + `
    "all analyses: participants were "

$newText = "We adjusted for the stratification variable (CHMP 2015) and the " + `
    "variables used to constrain randomization (Li 2017) as fixed effects in " + `
    "all analyses, using individual- rather than cluster-level measurements " + `
    "where possible. We followed the intention-to-treat principle for all " + `
    "analyses: participants were "

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false,
                         $true, 1, $false, $newText, 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Update the six odds-ratio summary statistics in the results table
#    (row "F": Odds Ratio, Robust Std. Err., z, P>|z|, [95% CI] low/high).
#    Cell.Range.Text is used instead of a document-wide Find/Replace so the
#    surrounding cell/paragraph formatting is left completely untouched.
# ---------------------------------------------------------------------------
$t = $d.Tables(1)
$t.Cell(5, 2).Range.Text = "0.33"
$t.Cell(5, 3).Range.Text = "0.13"
$t.Cell(5, 4).Range.Text = "2.55"
$t.Cell(5, 5).Range.Text = "0.01"
$t.Cell(5, 6).Range.Text = "0.08"
$t.Cell(5, 7).Range.Text = "0.58"
